# Update the code-review tracker sheet:
#  - row 20's "Вып." status (C20) moves from 3 to 1
#  - row 20's "Примечание" reply (E20) is replaced with the new reply text
#  - cosmetic: the selection/scroll position and row 20's height/column E
#    width are nudged to match where the author ended up after the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- core data edit -------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = "Прихуяриль параллельный парсинг по каждому сайту. Хз это ли имелось ввиду. Асинхронно вызывать GetHtmlDocumentFromWeb смысла особого нет, так как сразу же используется результат данного метода. Один хуй ждать придется"

# --- cosmetic follow-up (row height grew because the reply is longer) -
$ws.Rows(20).RowHeight = 75
$ws.Columns("E").ColumnWidth = 55.85546875

# --- cosmetic: cursor/scroll ended on D19 after the edit --------------
$ws.Range("D19").Select()
